$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("Metadata")

# --- Data sheet: update/extend the Fecha/Valor series ---
# Rows 2-40 now cover years 2023 down to 1985 (2 new recent years inserted
# at the top, plus 22 new older years appended at the bottom; every row's
# Valor also refreshed to match the latest published series).
    $ws1.Cells.Item(2,1).NumberFormat = "@"
    $ws1.Cells.Item(2,1).Value = "2023"
    $ws1.Cells.Item(2,2).Value = 50.8
    $ws1.Cells.Item(3,1).NumberFormat = "@"
    $ws1.Cells.Item(3,1).Value = "2022"
    $ws1.Cells.Item(3,2).Value = 50.9
    $ws1.Cells.Item(4,1).NumberFormat = "@"
    $ws1.Cells.Item(4,1).Value = "2021"
    $ws1.Cells.Item(4,2).Value = 51
    $ws1.Cells.Item(5,1).NumberFormat = "@"
    $ws1.Cells.Item(5,1).Value = "2020"
    $ws1.Cells.Item(5,2).Value = 52.6
    $ws1.Cells.Item(6,1).NumberFormat = "@"
    $ws1.Cells.Item(6,1).Value = "2019"
    $ws1.Cells.Item(6,2).Value = 50.6
    $ws1.Cells.Item(7,1).NumberFormat = "@"
    $ws1.Cells.Item(7,1).Value = "2018"
    $ws1.Cells.Item(7,2).Value = 50.8
    $ws1.Cells.Item(8,1).NumberFormat = "@"
    $ws1.Cells.Item(8,1).Value = "2017"
    $ws1.Cells.Item(8,2).Value = 50.5
    $ws1.Cells.Item(9,1).NumberFormat = "@"
    $ws1.Cells.Item(9,1).Value = "2016"
    $ws1.Cells.Item(9,2).Value = 50.3
    $ws1.Cells.Item(10,1).NumberFormat = "@"
    $ws1.Cells.Item(10,1).Value = "2015"
    $ws1.Cells.Item(10,2).Value = 50.8
    $ws1.Cells.Item(11,1).NumberFormat = "@"
    $ws1.Cells.Item(11,1).Value = "2014"
    $ws1.Cells.Item(11,2).Value = 50.2
    $ws1.Cells.Item(12,1).NumberFormat = "@"
    $ws1.Cells.Item(12,1).Value = "2013"
    $ws1.Cells.Item(12,2).Value = 50.2
    $ws1.Cells.Item(13,1).NumberFormat = "@"
    $ws1.Cells.Item(13,1).Value = "2012"
    $ws1.Cells.Item(13,2).Value = 50.3
    $ws1.Cells.Item(14,1).NumberFormat = "@"
    $ws1.Cells.Item(14,1).Value = "2011"
    $ws1.Cells.Item(14,2).Value = 51.1
    $ws1.Cells.Item(15,1).NumberFormat = "@"
    $ws1.Cells.Item(15,1).Value = "2010"
    $ws1.Cells.Item(15,2).Value = 51.6
    $ws1.Cells.Item(16,1).NumberFormat = "@"
    $ws1.Cells.Item(16,1).Value = "2009"
    $ws1.Cells.Item(16,2).Value = 50.1
    $ws1.Cells.Item(17,1).NumberFormat = "@"
    $ws1.Cells.Item(17,1).Value = "2008"
    $ws1.Cells.Item(17,2).Value = 50.1
    $ws1.Cells.Item(18,1).NumberFormat = "@"
    $ws1.Cells.Item(18,1).Value = "2007"
    $ws1.Cells.Item(18,2).Value = 53.7
    $ws1.Cells.Item(19,1).NumberFormat = "@"
    $ws1.Cells.Item(19,1).Value = "2006"
    $ws1.Cells.Item(19,2).Value = 55.5
    $ws1.Cells.Item(20,1).NumberFormat = "@"
    $ws1.Cells.Item(20,1).Value = "2005"
    $ws1.Cells.Item(20,2).Value = 56.8
    $ws1.Cells.Item(21,1).NumberFormat = "@"
    $ws1.Cells.Item(21,1).Value = "2004"
    $ws1.Cells.Item(21,2).Value = 60
    $ws1.Cells.Item(22,1).NumberFormat = "@"
    $ws1.Cells.Item(22,1).Value = "2003"
    $ws1.Cells.Item(22,2).Value = 61.5
    $ws1.Cells.Item(23,1).NumberFormat = "@"
    $ws1.Cells.Item(23,1).Value = "2002"
    $ws1.Cells.Item(23,2).Value = 63.6
    $ws1.Cells.Item(24,1).NumberFormat = "@"
    $ws1.Cells.Item(24,1).Value = "2001"
    $ws1.Cells.Item(24,2).Value = 62.8
    $ws1.Cells.Item(25,1).NumberFormat = "@"
    $ws1.Cells.Item(25,1).Value = "2000"
    $ws1.Cells.Item(25,2).Value = 64
    $ws1.Cells.Item(26,1).NumberFormat = "@"
    $ws1.Cells.Item(26,1).Value = "1999"
    $ws1.Cells.Item(26,2).Value = 64
    $ws1.Cells.Item(27,1).NumberFormat = "@"
    $ws1.Cells.Item(27,1).Value = "1998"
    $ws1.Cells.Item(27,2).Value = 61.8
    $ws1.Cells.Item(28,1).NumberFormat = "@"
    $ws1.Cells.Item(28,1).Value = "1997"
    $ws1.Cells.Item(28,2).Value = 64.5
    $ws1.Cells.Item(29,1).NumberFormat = "@"
    $ws1.Cells.Item(29,1).Value = "1996"
    $ws1.Cells.Item(29,2).Value = 64.3
    $ws1.Cells.Item(30,1).NumberFormat = "@"
    $ws1.Cells.Item(30,1).Value = "1995"
    $ws1.Cells.Item(30,2).Value = 68.1
    $ws1.Cells.Item(31,1).NumberFormat = "@"
    $ws1.Cells.Item(31,1).Value = "1994"
    $ws1.Cells.Item(31,2).Value = 67.3
    $ws1.Cells.Item(32,1).NumberFormat = "@"
    $ws1.Cells.Item(32,1).Value = "1993"
    $ws1.Cells.Item(32,2).Value = 67.5
    $ws1.Cells.Item(33,1).NumberFormat = "@"
    $ws1.Cells.Item(33,1).Value = "1992"
    $ws1.Cells.Item(33,2).Value = 67.6
    $ws1.Cells.Item(34,1).NumberFormat = "@"
    $ws1.Cells.Item(34,1).Value = "1991"
    $ws1.Cells.Item(34,2).Value = 66.1
    $ws1.Cells.Item(35,1).NumberFormat = "@"
    $ws1.Cells.Item(35,1).Value = "1990"
    $ws1.Cells.Item(35,2).Value = 65.1
    $ws1.Cells.Item(36,1).NumberFormat = "@"
    $ws1.Cells.Item(36,1).Value = "1989"
    $ws1.Cells.Item(36,2).Value = 65
    $ws1.Cells.Item(37,1).NumberFormat = "@"
    $ws1.Cells.Item(37,1).Value = "1988"
    $ws1.Cells.Item(37,2).Value = 63.3
    $ws1.Cells.Item(38,1).NumberFormat = "@"
    $ws1.Cells.Item(38,1).Value = "1987"
    $ws1.Cells.Item(38,2).Value = 62.7
    $ws1.Cells.Item(39,1).NumberFormat = "@"
    $ws1.Cells.Item(39,1).Value = "1986"
    $ws1.Cells.Item(39,2).Value = 61.7
    $ws1.Cells.Item(40,1).NumberFormat = "@"
    $ws1.Cells.Item(40,1).Value = "1985"
    $ws1.Cells.Item(40,2).Value = 64.1

# --- Metadata sheet: add "actualizacion" row before "cita", shifting the
#     trailing rows down by one ---
$ws2.Cells.Item(1,1).Value = " "
$ws2.Cells.Item(9,1).Value = "actualizacion"
$ws2.Cells.Item(9,2).Value = "Julio 2025"
$ws2.Cells.Item(10,1).Value = "cita"
$ws2.Cells.Item(10,2).Value = "UMAD con base en DINEM - MIDES hasta 2018, a partir de 2019 MIDES-MEF-OPP"
$ws2.Cells.Item(11,1).Value = "Mirador DESCA - UMAD/FCS – INDDHH"
$ws2.Cells.Item(11,2).Value = " "
